$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values - repulled data / recalculated mean
$ws.Range("F4").Value = -11
$ws.Range("F5").Value = -9
$ws.Range("F6").Value = 4
$ws.Range("F8").Value = -6
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -3
